# Better detection of numerical data types. Not only System.Double.
#
# The "Data Types" sample sheet gains three new example rows right after
# the existing "TimeSpan:" row (and before "Explicit Text:"):
#   Decimal Number:  123.45   (numeric)
#   Float Number:    123.45   (numeric)
#   Double Number:   123.45   (numeric)
#
# This pushes every row from the old row 9 ("Explicit Text:") onward down
# by four rows (three new data rows plus the existing blank separator row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four new rows starting at row 9 - this shifts the "Explicit Text:"
# section (and everything below it) down, preserving the blank separator
# row that used to sit at row 8 and keeping it at the same relative offset.
$ws.Rows("9:12").Insert()

$ws.Range("B9").Value = "Decimal Number:"
$ws.Range("C9").Value = 123.45

$ws.Range("B10").Value = "Float Number:"
$ws.Range("C10").Value = 123.45

$ws.Range("B11").Value = "Double Number:"
$ws.Range("C11").Value = 123.45
